# Pokazuj średnią netto z całego roku
# Adds a new "NettoSr" (average net) row under the existing "Netto" totals
# row (row 27) on every monthly sheet, showing Suma(Netto)/12 in column E.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A28").Value = "NettoSr"
    $ws.Range("B28").Value = "netto srednio"
    $ws.Range("C28").Value = "Suma(Netto)/12"
    $ws.Range("D28").Value = "zł"

    $ws.Range("E28").Formula = "=ROUND(SUM(E27:P27)/12,2)"
    $ws.Range("E28").NumberFormat = "#,##0.00"
}
